$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = "ba"
$ws.Range("J14").Value = "Appreciation"
$ws.Range("I25").Value = "b"
$ws.Range("J25").Value = "Acknowledge (Backchannel)"
$ws.Range("I33").Value = "b"
$ws.Range("J33").Value = "Acknowledge (Backchannel)"
$ws.Range("I41").Value = "sd"
$ws.Range("J41").Value = "Statement-non-opinion"
$ws.Range("I50").Value = "sv"
$ws.Range("J50").Value = "Statement-opinion"
$ws.Range("I61").Value = "aa"
$ws.Range("J61").Value = "Agree/Accept"
$ws.Range("I62").Value = "sd"
$ws.Range("J62").Value = "Statement-non-opinion"
$ws.Range("I64").Value = "ba"
$ws.Range("J64").Value = "Appreciation"
$ws.Range("I68").Value = "aa"
$ws.Range("J68").Value = "Agree/Accept"
$ws.Range("I82").Value = "ba"
$ws.Range("J82").Value = "Appreciation"
$ws.Range("I99").Value = "sd"
$ws.Range("J99").Value = "Statement-non-opinion"
$ws.Range("I111").Value = "sv"
$ws.Range("J111").Value = "Statement-opinion"
$ws.Range("I112").Value = "ba"
$ws.Range("J112").Value = "Appreciation"
$ws.Range("I121").Value = "ba"
$ws.Range("J121").Value = "Appreciation"
$ws.Range("I129").Value = "ba"
$ws.Range("J129").Value = "Appreciation"
$ws.Range("I134").Value = "sd"
$ws.Range("J134").Value = "Statement-non-opinion"
$ws.Range("I139").Value = "aa"
$ws.Range("J139").Value = "Agree/Accept"
$ws.Range("I140").Value = "aa"
$ws.Range("J140").Value = "Agree/Accept"
$ws.Range("I152").Value = "ba"
$ws.Range("J152").Value = "Appreciation"
$ws.Range("I164").Value = "%"
$ws.Range("J164").Value = "Uninterpretable"
$ws.Range("I172").Value = "sd"
$ws.Range("J172").Value = "Statement-non-opinion"
$ws.Range("I180").Value = "sd"
$ws.Range("J180").Value = "Statement-non-opinion"
$ws.Range("I188").Value = "%"
$ws.Range("J188").Value = "Uninterpretable"
$ws.Range("I193").Value = "sd"
$ws.Range("J193").Value = "Statement-non-opinion"
$ws.Range("I194").Value = "sd"
$ws.Range("J194").Value = "Statement-non-opinion"
$ws.Range("I199").Value = "sd"
$ws.Range("J199").Value = "Statement-non-opinion"
$ws.Range("I205").Value = "sd"
$ws.Range("J205").Value = "Statement-non-opinion"
$ws.Range("I206").Value = "sd"
$ws.Range("J206").Value = "Statement-non-opinion"
$ws.Range("I222").Value = "qy"
$ws.Range("J222").Value = "Yes-No-Question"
$ws.Range("I230").Value = "sd"
$ws.Range("J230").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "%"
$ws.Range("J243").Value = "Uninterpretable"
$ws.Range("I244").Value = "%"
$ws.Range("J244").Value = "Uninterpretable"
$ws.Range("I245").Value = "sv"
$ws.Range("J245").Value = "Statement-opinion"
$ws.Range("I250").Value = "sv"
$ws.Range("J250").Value = "Statement-opinion"
$ws.Range("I268").Value = "sd"
$ws.Range("J268").Value = "Statement-non-opinion"
$ws.Range("I286").Value = "sv"
$ws.Range("J286").Value = "Statement-opinion"
$ws.Range("I307").Value = "sv"
$ws.Range("J307").Value = "Statement-opinion"
$ws.Range("I312").Value = "sd"
$ws.Range("J312").Value = "Statement-non-opinion"
$ws.Range("I314").Value = "sd"
$ws.Range("J314").Value = "Statement-non-opinion"
$ws.Range("I317").Value = "aa"
$ws.Range("J317").Value = "Agree/Accept"
$ws.Range("I329").Value = "sd"
$ws.Range("J329").Value = "Statement-non-opinion"
